# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) table for rows 16-22
# gets its period order reversed (1908..1902 -> 1902..1908), carrying the
# matching Valor Mora amount along with each period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1902"
$ws.Range("F16").Value = 33125

$ws.Range("E17").Value = "1903"
$ws.Range("F17").Value = 33125

$ws.Range("E18").Value = "1904"
$ws.Range("F18").Value = 33125

$ws.Range("E19").Value = "1905"
$ws.Range("F19").Value = 33125

$ws.Range("E20").Value = "1906"
$ws.Range("F20").Value = 33125

$ws.Range("E21").Value = "1907"
$ws.Range("F21").Value = 33125

$ws.Range("E22").Value = "1908"
$ws.Range("F22").Value = 20979
